$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1357.1428
$ws.Range("I70").Value = 1166.6666
$ws.Range("K70").Value = 3499.9998
$ws.Range("M70").Value = -3229.9998
$ws.Range("H73").Value = 1357.1428
$ws.Range("I73").Value = 1166.6666
$ws.Range("K73").Value = 3499.9998
$ws.Range("M73").Value = -2563.9998
$ws.Range("H80").Value = 13333880
$ws.Range("I80").Value = 27778412
$ws.Range("J80").Value = 465.30768
$ws.Range("K80").Value = 83335236
$ws.Range("L80").Value = 1395.92304
$ws.Range("M80").Value = -83334238
$ws.Range("N80").Value = -3391.92304
$ws.Range("H83").Value = 13333880
$ws.Range("I83").Value = 27778412
$ws.Range("J83").Value = 465.30768
$ws.Range("K83").Value = 250005708
$ws.Range("L83").Value = 4187.76912
$ws.Range("M83").Value = -250000716
$ws.Range("N83").Value = -14171.76912
$ws.Range("H116").Value = 6212.8423
$ws.Range("J116").Value = 7470.0835
$ws.Range("L116").Value = 7470.0835
$ws.Range("N116").Value = -14354.0835
$ws.Range("H120").Value = 49726
$ws.Range("J120").Value = 49726
$ws.Range("L120").Value = 49726
$ws.Range("N120").Value = -59402
$ws.Range("H132").Value = 28510.344
$ws.Range("I132").Value = 4112.8213
$ws.Range("J132").Value = 126100.43
$ws.Range("K132").Value = 12338.4639
$ws.Range("L132").Value = 378301.29
$ws.Range("M132").Value = -9808.463899999999
$ws.Range("N132").Value = -383361.29
$ws.Range("H135").Value = 16667887
$ws.Range("I135").Value = 1089.909
$ws.Range("K135").Value = 9809.181
$ws.Range("M135").Value = -7274.181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35790.48
$ws.Range("I32").Value = 36604.566
$ws.Range("J32").Value = 29820.5
$ws.Range("K32").Value = 36604.566
$ws.Range("L32").Value = 29820.5
$ws.Range("M32").Value = -36317.566
$ws.Range("N32").Value = -30394.5
$ws.Range("H80").Value = 54996
$ws.Range("J80").Value = 54996
$ws.Range("L80").Value = 54996
$ws.Range("N80").Value = -56992
$ws.Range("H83").Value = 54996
$ws.Range("J83").Value = 54996
$ws.Range("L83").Value = 164988
$ws.Range("N83").Value = -174972
$ws.Range("H113").Value = 40196
$ws.Range("J113").Value = 40196
$ws.Range("L113").Value = 40196
$ws.Range("N113").Value = -48874
$ws.Range("H117").Value = 47311.75
$ws.Range("J117").Value = 47311.75
$ws.Range("L117").Value = 47311.75
$ws.Range("N117").Value = -56489.75
$ws.Range("H118").Value = 49803
$ws.Range("J118").Value = 49803
$ws.Range("L118").Value = 49803
$ws.Range("N118").Value = -53117

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 42496
$ws.Range("J112").Value = 42496
$ws.Range("L112").Value = 42496
$ws.Range("N112").Value = -45450
$ws.Range("H117").Value = 49914
$ws.Range("J117").Value = 49914
$ws.Range("L117").Value = 49914
$ws.Range("N117").Value = -59092
$ws.Range("H126").Value = 50780
$ws.Range("J126").Value = 50780
$ws.Range("L126").Value = 50780
$ws.Range("N126").Value = -60660
$ws.Range("H132").Value = 60780
$ws.Range("J132").Value = 60780
$ws.Range("L132").Value = 60780
$ws.Range("N132").Value = -70900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 36025.332
$ws.Range("J110").Value = 36025.332
$ws.Range("L110").Value = 36025.332
$ws.Range("N110").Value = -44205.332
$ws.Range("H111").Value = 46994.668
$ws.Range("J111").Value = 46994.668
$ws.Range("L111").Value = 46994.668
$ws.Range("N111").Value = -55174.668
$ws.Range("H118").Value = 44742
$ws.Range("J118").Value = 44742
$ws.Range("L118").Value = 44742
$ws.Range("N118").Value = -48056
$ws.Range("H122").Value = 240840
$ws.Range("I122").Value = 300750
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 902250
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -899800
$ws.Range("N122").Value = -8500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3563.2327
$ws.Range("I131").Value = 20697.6
$ws.Range("J131").Value = 1308.7106
$ws.Range("K131").Value = 62092.8
$ws.Range("L131").Value = 3926.1318
$ws.Range("M131").Value = -57052.8
$ws.Range("N131").Value = -14006.1318

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 221106.73
$ws.Range("I80").Value = 337067
$ws.Range("J80").Value = 3681.25
$ws.Range("K80").Value = 337067
$ws.Range("L80").Value = 3681.25
$ws.Range("M80").Value = -336069
$ws.Range("N80").Value = -5677.25
$ws.Range("H83").Value = 221106.73
$ws.Range("I83").Value = 337067
$ws.Range("J83").Value = 3681.25
$ws.Range("K83").Value = 1685335
$ws.Range("L83").Value = 18406.25
$ws.Range("M83").Value = -1680343
$ws.Range("N83").Value = -28390.25
$ws.Range("H102").Value = 3169.6667
$ws.Range("I102").Value = 3000.8
$ws.Range("K102").Value = 3000.8
$ws.Range("M102").Value = -1378.8
$ws.Range("H110").Value = 46718.8
$ws.Range("J110").Value = 46718.8
$ws.Range("L110").Value = 46718.8
$ws.Range("N110").Value = -54898.8
$ws.Range("H116").Value = 39997.332
$ws.Range("J116").Value = 39997.332
$ws.Range("L116").Value = 39997.332
$ws.Range("N116").Value = -49175.332
$ws.Range("H119").Value = 33202
$ws.Range("J119").Value = 33202
$ws.Range("L119").Value = 33202
$ws.Range("N119").Value = -42878
$ws.Range("H133").Value = 29904.21
$ws.Range("J133").Value = 29904.21
$ws.Range("L133").Value = 29904.21
$ws.Range("N133").Value = -40024.21
$ws.Range("H135").Value = 45239.8
$ws.Range("J135").Value = 45239.8
$ws.Range("L135").Value = 45239.8
$ws.Range("N135").Value = -55379.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2323.9412
$ws.Range("I7").Value = 1932.909
$ws.Range("J7").Value = 3040.8333
$ws.Range("K7").Value = 1932.909
$ws.Range("L7").Value = 3040.8333
$ws.Range("M7").Value = -1820.909
$ws.Range("N7").Value = -3264.8333
$ws.Range("H48").Value = 10000
$ws.Range("J48").Value = 10000
$ws.Range("L48").Value = 10000
$ws.Range("N48").Value = -11322
$ws.Range("H81").Value = 42181
$ws.Range("J81").Value = 42181
$ws.Range("L81").Value = 42181
$ws.Range("N81").Value = -44177
$ws.Range("H82").Value = 6945451.5
$ws.Range("I82").Value = 1154.2858
$ws.Range("K82").Value = 1154.2858
$ws.Range("M82").Value = -793.2858000000001
$ws.Range("H84").Value = 42181
$ws.Range("J84").Value = 42181
$ws.Range("L84").Value = 126543
$ws.Range("N84").Value = -136527
$ws.Range("H85").Value = 6945451.5
$ws.Range("I85").Value = 1154.2858
$ws.Range("K85").Value = 1154.2858
$ws.Range("M85").Value = 93.71419999999989
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H118").Value = 41142
$ws.Range("J118").Value = 41142
$ws.Range("L118").Value = 41142
$ws.Range("N118").Value = -44456
$ws.Range("H121").Value = 42280
$ws.Range("J121").Value = 42280
$ws.Range("L121").Value = 42280
$ws.Range("N121").Value = -45774
$ws.Range("H126").Value = 2323.9412
$ws.Range("I126").Value = 1932.909
$ws.Range("J126").Value = 3040.8333
$ws.Range("K126").Value = 5798.727000000001
$ws.Range("L126").Value = 9122.499899999999
$ws.Range("M126").Value = -3328.727000000001
$ws.Range("N126").Value = -14062.4999
$ws.Range("H127").Value = 49416.5
$ws.Range("J127").Value = 49416.5
$ws.Range("L127").Value = 49416.5
$ws.Range("N127").Value = -59336.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 61446.547
$ws.Range("J46").Value = 61446.547
$ws.Range("L46").Value = 61446.547
$ws.Range("N46").Value = -61908.547
$ws.Range("H81").Value = 1852.7333
$ws.Range("I81").Value = 1398.7142
$ws.Range("J81").Value = 2250
$ws.Range("K81").Value = 2797.4284
$ws.Range("L81").Value = 4500
$ws.Range("M81").Value = -1736.4284
$ws.Range("N81").Value = -6622
$ws.Range("H84").Value = 1852.7333
$ws.Range("I84").Value = 1398.7142
$ws.Range("J84").Value = 2250
$ws.Range("K84").Value = 13987.142
$ws.Range("L84").Value = 22500
$ws.Range("M84").Value = -8683.142
$ws.Range("N84").Value = -33108
$ws.Range("H108").Value = 42000
$ws.Range("J108").Value = 42000
$ws.Range("L108").Value = 42000
$ws.Range("N108").Value = -49680
$ws.Range("H110").Value = 48644
$ws.Range("J110").Value = 48644
$ws.Range("L110").Value = 48644
$ws.Range("N110").Value = -56824
$ws.Range("H112").Value = 33687
$ws.Range("J112").Value = 33687
$ws.Range("L112").Value = 33687
$ws.Range("N112").Value = -36641
$ws.Range("H117").Value = 41996
$ws.Range("J117").Value = 41996
$ws.Range("L117").Value = 41996
$ws.Range("N117").Value = -51174
$ws.Range("H119").Value = 48698
$ws.Range("J119").Value = 48698
$ws.Range("L119").Value = 48698
$ws.Range("N119").Value = -58374
$ws.Range("H120").Value = 42104
$ws.Range("J120").Value = 42104
$ws.Range("L120").Value = 42104
$ws.Range("N120").Value = -51780
$ws.Range("H134").Value = 61446.547
$ws.Range("J134").Value = 61446.547
$ws.Range("L134").Value = 184339.641
$ws.Range("N134").Value = -189409.641
